$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.033721973275866
$ws.Range("D2").Value = 1.042649370212068
$ws.Range("E2").Value = 1.051610385697098
$ws.Range("F2").Value = 1.056756728469804
$ws.Range("I2").Value = 1.038531091474237
$ws.Range("J2").Value = 1.038845142350186
$ws.Range("K2").Value = 1.045425461146601
$ws.Range("L2").Value = 1.054361410229219
$ws.Range("M2").Value = 1.059493565495453
$ws.Range("N2").Value = 1.040320421609527
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03465313759208
$ws.Range("D3").Value = 1.04335971836418
$ws.Range("E3").Value = 1.052504612184161
$ws.Range("F3").Value = 1.057626584507984
$ws.Range("I3").Value = 1.038717252178612
$ws.Range("J3").Value = 1.039419231887186
$ws.Range("K3").Value = 1.04594694719802
$ws.Range("L3").Value = 1.055068107682653
$ws.Range("M3").Value = 1.060176980543706
$ws.Range("N3").Value = 1.040895326419518
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.035255964416308
$ws.Range("D4").Value = 1.043819356766045
$ws.Range("E4").Value = 1.053084226840176
$ws.Range("F4").Value = 1.058189971051008
$ws.Range("I4").Value = 1.038836072506762
$ws.Range("J4").Value = 1.039790388737462
$ws.Range("K4").Value = 1.046283706624339
$ws.Range("L4").Value = 1.055525741405233
$ws.Range("M4").Value = 1.060619085110495
$ws.Range("N4").Value = 1.041267010355094
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.035509463894888
$ws.Range("D5").Value = 1.044012586219315
$ws.Range("E5").Value = 1.053328132383304
$ws.Range("F5").Value = 1.058426944438578
$ws.Range("I5").Value = 1.038885631681834
$ws.Range("J5").Value = 1.039946346059297
$ws.Range("K5").Value = 1.046425117097636
$ws.Range("L5").Value = 1.055718213885554
$ws.Range("M5").Value = 1.060804918051059
$ws.Range("N5").Value = 1.041423189154214
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.035552031710929
$ws.Range("D6").Value = 1.044045030102976
$ws.Range("E6").Value = 1.0533690989646
$ws.Range("F6").Value = 1.058466740637833
$ws.Range("I6").Value = 1.038893929824892
$ws.Range("J6").Value = 1.039972527433946
$ws.Range("K6").Value = 1.046448850935805
$ws.Range("L6").Value = 1.055750535706781
$ws.Range("M6").Value = 1.060836118549567
$ws.Range("N6").Value = 1.041449407709419
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.035259351410743
$ws.Range("D7").Value = 1.043821938719281
$ws.Range("E7").Value = 1.053087484994404
$ws.Range("F7").Value = 1.058193137009123
$ws.Range("I7").Value = 1.038836736263594
$ws.Range("J7").Value = 1.039792472950736
$ws.Range("K7").Value = 1.046285596800657
$ws.Range("L7").Value = 1.0555283129064
$ws.Range("M7").Value = 1.060621568330323
$ws.Range("N7").Value = 1.04126909752819
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.034036601617397
$ws.Range("D8").Value = 1.042889435703365
$ws.Range("E8").Value = 1.051912387744539
$ws.Range("F8").Value = 1.057050589397378
$ws.Range("I8").Value = 1.038594344097919
$ws.Range("J8").Value = 1.039039223721246
$ws.Range("K8").Value = 1.045601839279557
$ws.Range("L8").Value = 1.054600168027766
$ws.Range("M8").Value = 1.05972455092952
$ws.Range("N8").Value = 1.040514778598395
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031884313069602
$ws.Range("D9").Value = 1.04124627841246
$ws.Range("E9").Value = 1.04984936901922
$ws.Range("F9").Value = 1.055041414145639
$ws.Range("I9").Value = 1.038154704219703
$ws.Range("J9").Value = 1.037709515561607
$ws.Range("K9").Value = 1.044391836599019
$ws.Range("L9").Value = 1.052967420522561
$ws.Range("M9").Value = 1.058143102405258
$ws.Range("N9").Value = 1.039183182100649
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.030451095928143
$ws.Range("D10").Value = 1.040150949297185
$ws.Range("E10").Value = 1.048479259540487
$ws.Range("F10").Value = 1.053704839083325
$ws.Range("I10").Value = 1.037853239792202
$ws.Range("J10").Value = 1.036821498017306
$ws.Range("K10").Value = 1.043581782282701
$ws.Range("L10").Value = 1.051880855800657
$ws.Range("M10").Value = 1.057088350527776
$ws.Range("N10").Value = 1.038293903469579
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029830898675864
$ws.Range("D11").Value = 1.039676703747992
$ws.Range("E11").Value = 1.047887247981223
$ws.Range("F11").Value = 1.053126788353791
$ws.Range("I11").Value = 1.037720725588308
$ws.Range("J11").Value = 1.036436622273399
$ws.Range("K11").Value = 1.043230230453334
$ws.Range("L11").Value = 1.051410835795375
$ws.Range("M11").Value = 1.056631540042251
$ws.Range("N11").Value = 1.03790848115797
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029600590100707
$ws.Range("D12").Value = 1.039500555157505
$ws.Range("E12").Value = 1.047667538573367
$ws.Range("F12").Value = 1.052912180515022
$ws.Range("I12").Value = 1.037671207387552
$ws.Range("J12").Value = 1.036293609402768
$ws.Range("K12").Value = 1.043099530383339
$ws.Range("L12").Value = 1.051236321260379
$ws.Range("M12").Value = 1.056461847137824
$ws.Range("N12").Value = 1.037765265192659
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029649989328452
$ws.Range("D13").Value = 1.039538339273775
$ws.Range("E13").Value = 1.047714658362998
$ws.Range("F13").Value = 1.052958209806722
$ws.Range("I13").Value = 1.03768184261745
$ws.Range("J13").Value = 1.036324288526875
$ws.Range("K13").Value = 1.043127571317749
$ws.Range("L13").Value = 1.051273751950249
$ws.Range("M13").Value = 1.056498247416642
$ws.Range("N13").Value = 1.03779598788464
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029811860061795
$ws.Range("D14").Value = 1.039662143092898
$ws.Range("E14").Value = 1.04786908284744
$ws.Range("F14").Value = 1.053109046633026
$ws.Range("I14").Value = 1.037716638443574
$ws.Range("J14").Value = 1.036424801858227
$ws.Range("K14").Value = 1.043219429147165
$ws.Range("L14").Value = 1.051396408897047
$ws.Range("M14").Value = 1.056617513434786
$ws.Range("N14").Value = 1.037896643956452
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029911601980721
$ws.Range("D15").Value = 1.039738423725017
$ws.Range("E15").Value = 1.047964254096905
$ws.Range("F15").Value = 1.05320199626314
$ws.Range("I15").Value = 1.03773803802783
$ws.Range("J15").Value = 1.036486724459818
$ws.Range("K15").Value = 1.043276010174283
$ws.Range("L15").Value = 1.051471991444924
$ws.Range("M15").Value = 1.056690995460101
$ws.Range("N15").Value = 1.03795865449524
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.030492264782708
$ws.Range("D16").Value = 1.040182424360596
$ws.Range("E16").Value = 1.048518575966058
$ws.Range("F16").Value = 1.053743217180246
$ws.Range("I16").Value = 1.037861992703599
$ws.Range("J16").Value = 1.036847033484948
$ws.Range("K16").Value = 1.043605097004038
$ws.Range("L16").Value = 1.051912059483596
$ws.Range("M16").Value = 1.057118665621594
$ws.Range("N16").Value = 1.038319475200514
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030856605680375
$ws.Range("D17").Value = 1.040460945734692
$ws.Range("E17").Value = 1.048866624466368
$ws.Range("F17").Value = 1.054082897925677
$ws.Range("I17").Value = 1.03793921706846
$ws.Range("J17").Value = 1.03707295056184
$ws.Range("K17").Value = 1.043811312932537
$ws.Range("L17").Value = 1.052188229175053
$ws.Range("M17").Value = 1.05738690701821
$ws.Range("N17").Value = 1.038545713105568
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.031069157633808
$ws.Range("D18").Value = 1.040623406224887
$ws.Range("E18").Value = 1.049069756245129
$ws.Range("F18").Value = 1.054281094922788
$ws.Range("I18").Value = 1.037984069829852
$ws.Range("J18").Value = 1.037204689537731
$ws.Range("K18").Value = 1.043931518613793
$ws.Range("L18").Value = 1.05234935951769
$ws.Range("M18").Value = 1.057543358342807
$ws.Range("N18").Value = 1.038677639165916
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.031141638786789
$ws.Range("D19").Value = 1.040678801634981
$ws.Range("E19").Value = 1.049139039406066
$ws.Range("F19").Value = 1.054348686289415
$ws.Range("I19").Value = 1.037999331061494
$ws.Range("J19").Value = 1.037249603220937
$ws.Range("K19").Value = 1.043972492622767
$ws.Range("L19").Value = 1.052404308434754
$ws.Range("M19").Value = 1.057596702577994
$ws.Range("N19").Value = 1.038722616631703
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030817511405479
$ws.Range("D20").Value = 1.040431062635607
$ws.Range("E20").Value = 1.048829269632015
$ws.Range("F20").Value = 1.054046446458929
$ws.Range("I20").Value = 1.037930951366238
$ws.Range("J20").Value = 1.037048715372919
$ws.Range("K20").Value = 1.043789195828407
$ws.Range("L20").Value = 1.052158594092931
$ws.Range("M20").Value = 1.057358128212252
$ws.Range("N20").Value = 1.038521443499899
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.029764191493413
$ws.Range("D21").Value = 1.039625685740263
$ws.Range("E21").Value = 1.047823603427407
$ws.Range("F21").Value = 1.053064626010616
$ws.Range("I21").Value = 1.037706400120264
$ws.Range("J21").Value = 1.036395204639608
$ws.Range("K21").Value = 1.043192382556597
$ws.Range("L21").Value = 1.051360287490617
$ws.Range("M21").Value = 1.056582392917932
$ws.Range("N21").Value = 1.037867004706388
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.029102276784452
$ws.Range("D22").Value = 1.03911935639225
$ws.Range("E22").Value = 1.047192401456281
$ws.Range("F22").Value = 1.05244793020747
$ws.Range("I22").Value = 1.037563500563535
$ws.Range("J22").Value = 1.035984011007554
$ws.Range("K22").Value = 1.0428164597874
$ws.Range("L22").Value = 1.050858776609538
$ws.Range("M22").Value = 1.056094581889362
$ws.Range("N22").Value = 1.037455227132212
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029453136875822
$ws.Range("D23").Value = 1.039387766599077
$ws.Range("E23").Value = 1.047526908765044
$ws.Range("F23").Value = 1.052774793647619
$ws.Range("I23").Value = 1.037639416689265
$ws.Range("J23").Value = 1.036202021174082
$ws.Range("K23").Value = 1.04301580790529
$ws.Range("L23").Value = 1.051124597149625
$ws.Range("M23").Value = 1.056353186475776
$ws.Range("N23").Value = 1.037673546898182
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.030835176307918
$ws.Range("D24").Value = 1.040444565510537
$ws.Range("E24").Value = 1.048846148300097
$ws.Range("F24").Value = 1.054062917100954
$ws.Range("I24").Value = 1.037934686871193
$ws.Range("J24").Value = 1.037059666318685
$ws.Range("K24").Value = 1.043799189832458
$ws.Range("L24").Value = 1.052171984771273
$ws.Range("M24").Value = 1.057371132144883
$ws.Range("N24").Value = 1.038532409997264
$ws.Range("B25").Value = 1.019999999999999
$ws.Range("C25").Value = 1.03244044624955
$ws.Range("D25").Value = 1.041671061116533
$ws.Range("E25").Value = 1.050381792528144
$ws.Range("F25").Value = 1.055560334064334
$ws.Range("I25").Value = 1.038269840243844
$ws.Range("J25").Value = 1.038053553720902
$ws.Range("K25").Value = 1.039527708833593
$ws.Range("L25").Value = 1.05338918965976
$ws.Range("M25").Value = 1.058552030260016
$ws.Range("N25").Value = 1.039527708833593
